# Remove the duplicated/garbled "貝賣" property-source string: fix the
# "土地" (land) sheet's mis-typed "033貝賣" register_reason to "買賣", and
# fix the "建物" (building) sheet's "貝賣" register_reason to the same
# "買賣" text, so the two cells share one (correct) string and the stray
# duplicate entry disappears from the shared-strings table.

$wb = $excel.ActiveWorkbook

$wsLand = $wb.Worksheets.Item("土地")
$wsLand.Range("G2").Value = "買賣"

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("G2").Value = "買賣"
